# Updated 7 test cases for DC Unit scripts with new loading method details
#
# Adds a new "DC Unit Loading Details" mini-table (Name / Current (DC Units) /
# Current (worst case)) to the "Add Devices Loop A", "Add Devices Loop B" and
# "Panel LED" sheets, and touches the page setup (portrait) on the three
# "Delete Devices" / "Panel LED" sheets.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Add Devices Loop A (column C, rows 1-3)
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("Add Devices Loop A")

$wsA.Range("C1").Value = "DC Unit Loading Details Name"
$null = $wsA.Range("A7").Copy()
$null = $wsA.Range("C1").PasteSpecial($xlPasteFormats)

$wsA.Range("C2").Value = "Current (DC Units)"
$null = $wsA.Range("A8").Copy()
$null = $wsA.Range("C2").PasteSpecial($xlPasteFormats)

$wsA.Range("C3").Value = "Current (worst case)"
$null = $wsA.Range("A8").Copy()
$null = $wsA.Range("C3").PasteSpecial($xlPasteFormats)

$null = $wsA.Range("C1:C3").Select()

# ---------------------------------------------------------------------------
# Panel LED (column F, rows 1-3)
# ---------------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("Panel LED")

$wsP.Range("F1").Value = "DC Unit Loading Details Name"
$null = $wsP.Range("A7").Copy()
$null = $wsP.Range("F1").PasteSpecial($xlPasteFormats)

$wsP.Range("F2").Value = "Current (DC Units)"
$null = $wsA.Range("A8").Copy()
$null = $wsP.Range("F2").PasteSpecial($xlPasteFormats)

$wsP.Range("F3").Value = "Current (worst case)"
$null = $wsA.Range("A8").Copy()
$null = $wsP.Range("F3").PasteSpecial($xlPasteFormats)

$wsP.Columns.Item(6).ColumnWidth = 25.5

$wsP.PageSetup.Orientation = 1

$null = $wsP.Range("F1:F3").Select()

# ---------------------------------------------------------------------------
# Delete Devices Loop A / B - page setup touched only
# ---------------------------------------------------------------------------
$wsDA = $wb.Worksheets.Item("Delete Devices Loop A")
$wsDA.PageSetup.Orientation = 1

$wsDB = $wb.Worksheets.Item("Delete Devices Loop B")
$wsDB.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Add Devices Loop B (column C, rows 1-3) - done last so this sheet ends up
# the active tab, matching the saved selection state of the workbook.
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("Add Devices Loop B")

$wsB.Range("C1").Value = "DC Unit Loading Details Name"
$null = $wsB.Range("A7").Copy()
$null = $wsB.Range("C1").PasteSpecial($xlPasteFormats)

$wsB.Range("C2").Value = "Current (DC Units)"
$null = $wsB.Range("A8").Copy()
$null = $wsB.Range("C2").PasteSpecial($xlPasteFormats)

$wsB.Range("C3").Value = "Current (worst case)"
$null = $wsB.Range("A8").Copy()
$null = $wsB.Range("C3").PasteSpecial($xlPasteFormats)

$null = $wsB.Range("C1:C3").Select()
